$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match formatting of the existing header cells (e.g. AC1): bold font,
# thin border all around, centered horizontally, top vertical alignment.
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data rows 2-40: Wins=86, Losses=76, Ties=0
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 86  # AD
    $ws.Cells.Item($r, 31).Value = 76  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
